# Update countries & provincias Spain
# Refresh of the "Pais" COVID dashboard sheet: new pull of case data, which
# also re-sorts a handful of rows (the sheet is kept sorted by "Casos
# totales" descending) and bumps the "Datos actualizados" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp: 21:15 -> 22:32
$ws.Range("A1").Value = "Datos actualizados a 20 de Junio de 2020 a las 22:32"

# Estados Unidos (row 4) - refreshed totals
$ws.Range("B4").Value = 2323055
$ws.Range("C4").Value = 25865
$ws.Range("D4").Value = 963563
$ws.Range("E4").Value = 1237616
$ws.Range("G4").Value = 469
$ws.Range("H4").Value = 121876

# India (row 7) - refreshed totals
$ws.Range("B7").Value = 411727
$ws.Range("C7").Value = 15915
$ws.Range("D7").Value = 228181
$ws.Range("E7").Value = 170269
$ws.Range("G7").Value = 307
$ws.Range("H7").Value = 13277

# Alemania (row 14) - refreshed totals
$ws.Range("B14").Value = 190992
$ws.Range("C14").Value = 332
$ws.Range("E14").Value = 7332

# Barein (row 50) - refreshed totals
$ws.Range("E50").Value = 5569
$ws.Range("G50").Value = 3
$ws.Range("H50").Value = 60

# Israel (row 51) - refreshed totals
$ws.Range("B51").Value = 20633
$ws.Range("C51").Value = 294
$ws.Range("E51").Value = 4742
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 305

# Costa de Marfil now overtakes Finlandia -> rows 74/75 swap order
$ws.Range("A74").Value = "Costa de Marfil"
$ws.Range("B74").Value = 7276
$ws.Range("C74").Value = 402
$ws.Range("D74").Value = 2992
$ws.Range("E74").Value = 4232
$ws.Range("G74").Value = 3
$ws.Range("H74").Value = 52

$ws.Range("A75").Value = "Finlandia"
$ws.Range("B75").Value = 7142
$ws.Range("C75").Value = 9
$ws.Range("D75").Value = 6200
$ws.Range("E75").Value = 616
$ws.Range("H75").Value = 326

# Guinea (row 82) - refreshed totals
$ws.Range("B82").Value = 4960
$ws.Range("C82").Value = 56
$ws.Range("D82").Value = 3630
$ws.Range("E82").Value = 1303

# Niger (row 126) - refreshed totals
$ws.Range("B126").Value = 1035
$ws.Range("C126").Value = 15
$ws.Range("D126").Value = 911
$ws.Range("E126").Value = 57

# Ruanda jumps ahead of San Marino, Santo Tome y Principe, Mozambique and
# Malta -> rows 140-144 each shift down one place
$ws.Range("A140").Value = "Ruanda"
$ws.Range("B140").Value = 702
$ws.Range("C140").Value = 41
$ws.Range("D140").Value = 357
$ws.Range("E140").Value = 343
$ws.Range("H140").Value = 2

$ws.Range("A141").Value = "San Marino"
$ws.Range("B141").Value = 696
$ws.Range("D141").Value = 610
$ws.Range("E141").Value = 44
$ws.Range("H141").Value = 42

$ws.Range("A142").Value = "Santo Tome y Principe"
$ws.Range("B142").Value = 693
$ws.Range("C142").Value = 0
$ws.Range("D142").Value = 199
$ws.Range("E142").Value = 482
$ws.Range("H142").Value = 12

$ws.Range("A143").Value = "Mozambique"
$ws.Range("B143").Value = 688
$ws.Range("C143").Value = 20
$ws.Range("D143").Value = 177
$ws.Range("E143").Value = 507
$ws.Range("H143").Value = 4

$ws.Range("A144").Value = "Malta"
$ws.Range("B144").Value = 664
$ws.Range("C144").Value = 1
$ws.Range("D144").Value = 616
$ws.Range("E144").Value = 39
$ws.Range("H144").Value = 9

# Togo (row 150) - refreshed totals
$ws.Range("B150").Value = 561
$ws.Range("C150").Value = 6
$ws.Range("D150").Value = 366
$ws.Range("E150").Value = 182

# Dominica / Fiyi tie-break swap (rows 202/203), identical case counts
$ws.Range("A202").Value = "Dominica"
$ws.Range("A203").Value = "Fiyi"

# Islas Turcas y Caicos / Santa Sede swap (rows 208/209)
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 1

$ws.Range("A209").Value = "Santa Sede"
$ws.Range("D209").Value = 12
$ws.Range("H209").Value = 0

# Islas Virgenes Britanicas / Papua Nueva Guinea swap (rows 213/214)
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0
